$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New departure rows for Friday, Jan 13 (appended after the existing 215 rows).
# Columns: Number, Date, Time, Flight, To, Short, Airline, Model, Aircraft ID, Status, Difference
$newRows = @(
    ,(215, "Friday, Jan 13", "5:35 AM", "LO3828", "Warsaw", "(WAW)", "LOT ", "E75S", "(SP-LIL)", "5:56 AM", "0 hours, 21 minutes")
    ,(216, "Friday, Jan 13", "6:00 AM", "FR6119", "Hamburg", "(HAM)", "Ryanair ", "B738", "(SP-RSW)", "6:22 AM", "0 hours, 22 minutes")
    ,(217, "Friday, Jan 13", "6:00 AM", "FR8782", "Cork", "(ORK)", "Ryanair ", "B738", "(SP-RSL)", "6:05 AM", "0 hours, 5 minutes")
    ,(218, "Friday, Jan 13", "6:00 AM", "SK1758", "Copenhagen", "(CPH)", "SAS ", "A20N", "(EI-SIU)", "6:13 AM", "0 hours, 13 minutes")
    ,(219, "Friday, Jan 13", "6:00 AM", "W61601", "London", "(LTN)", "Wizz Air ", "A320", "(HA-LYS)", "6:15 AM", "0 hours, 15 minutes")
    ,(220, "Friday, Jan 13", "6:05 AM", "FR6099", "Stockholm", "(ARN)", "Ryanair ", "B738", "(SP-RKM)", "6:17 AM", "0 hours, 12 minutes")
    ,(221, "Friday, Jan 13", "6:20 AM", "W61641", "Eindhoven", "(EIN)", "Wizz Air ", "A21N", "(9H-WAP)", "6:26 AM", "0 hours, 6 minutes")
    ,(222, "Friday, Jan 13", "6:30 AM", "FR6844", "Copenhagen", "(CPH)", "Ryanair ", "B738", "(SP-RKQ)", "6:43 AM", "0 hours, 13 minutes")
    ,(223, "Friday, Jan 13", "6:30 AM", "W61731", "Stockholm", "(NYO)", "Wizz Air ", "A321", "(HA-LTB)", "6:35 AM", "0 hours, 5 minutes")
    ,(224, "Friday, Jan 13", "6:35 AM", "FR6095", "Gothenburg", "(GOT)", "Ryanair ", "B738", "(SP-RSO)", "6:41 AM", "0 hours, 6 minutes")
    ,(225, "Friday, Jan 13", "6:55 AM", "LH1379", "Frankfurt", "(FRA)", "Lufthansa ", "CRJ9", "(D-ACKI)", "7:01 AM", "0 hours, 6 minutes")
    ,(226, "Friday, Jan 13", "7:00 AM", "LPR42", "Warsaw", "(WAW)", "Polish Medical Air Rescue ", "LJ75", "(SP-MXS)", "7:25 AM", "0 hours, 25 minutes")
    ,(227, "Friday, Jan 13", "7:55 AM", "W61745", "Bergen", "(BGO)", "Wizz Air ", "A321", "(HA-LXL)", "7:58 AM", "0 hours, 3 minutes")
    ,(228, "Friday, Jan 13", "8:00 AM", "FR4106", "Wroclaw", "(WRO)", "Ryanair ", "B738", "(SP-RKG)", "8:12 AM", "0 hours, 12 minutes")
    ,(229, "Friday, Jan 13", "8:50 AM", "LO3832", "Warsaw", "(WAW)", "LOT (Warmia Mazury Livery) ", "E195", "(SP-LNF)", "9:00 AM", "0 hours, 10 minutes")
    ,(230, "Friday, Jan 13", "9:30 AM", "FR3686", "Billund", "(BLL)", "Ryanair ", "B738", "(SP-RSW)", "9:52 AM", "0 hours, 22 minutes")
    ,(231, "Friday, Jan 13", "9:35 AM", "FR6082", "London", "(STN)", "Ryanair ", "B38M", "(EI-HGT)", "10:01 AM", "0 hours, 26 minutes")
    ,(232, "Friday, Jan 13", "9:40 AM", "FR3277", "Oslo", "(TRF)", "Ryanair ", "B738", "(SP-RKQ)", "9:49 AM", "0 hours, 9 minutes")
    ,(233, "Friday, Jan 13", "9:45 AM", "FR827", "Venice", "(VCE)", "Ryanair ", "B738", "(SP-RKM)", "9:59 AM", "0 hours, 14 minutes")
    ,(234, "Friday, Jan 13", "9:55 AM", "FR3279", "Milan", "(BGY)", "Ryanair ", "B738", "(SP-RSO)", "10:06 AM", "0 hours, 11 minutes")
    ,(235, "Friday, Jan 13", "10:50 AM", "W61743", "Oslo", "(TRF)", "Wizz Air ", "A320", "(HA-LYM)", "10:55 AM", "0 hours, 5 minutes")
    ,(236, "Friday, Jan 13", "11:45 AM", "UNKNOWN", "Tuzla", "(TZL)", "Wizz Air ", "A320", "(HA-LYH)", "11:51 AM", "0 hours, 6 minutes")
    ,(237, "Friday, Jan 13", "11:55 AM", "LO3836", "Warsaw", "(WAW)", "LOT ", "E170", "(SP-LDH)", "11:54 AM", "0 hours, -1 minutes")
    ,(238, "Friday, Jan 13", "12:00 PM", "SAR900", "Warsaw", "(WAW)", "SprintAir ", "SF34", "(SP-KPV)", "12:08 PM", "0 hours, 8 minutes")
    ,(239, "Friday, Jan 13", "12:35 PM", "W61631", "Paris", "(BVA)", "Wizz Air ", "A320", "(HA-LYO)", "12:36 PM", "0 hours, 1 minutes")
    ,(240, "Friday, Jan 13", "1:00 PM", "DY1031", "Bergen", "(BGO)", "Norwegian ", "B738", "(LN-ENM)", "1:02 PM", "0 hours, 2 minutes")
    ,(241, "Friday, Jan 13", "1:10 PM", "DY1053", "Oslo", "(OSL)", "Norwegian ", "B738", "(SE-RPG)", "1:18 PM", "0 hours, 8 minutes")
    ,(242, "Friday, Jan 13", "1:10 PM", "LH1643", "Munich", "(MUC)", "Lufthansa ", "CRJ9", "(D-ACNX)", "1:22 PM", "0 hours, 12 minutes")
    ,(243, "Friday, Jan 13", "1:20 PM", "FR3299", "Birmingham", "(BHX)", "Ryanair ", "B738", "(SP-RSL)", "1:28 PM", "0 hours, 8 minutes")
    ,(244, "Friday, Jan 13", "1:30 PM", "FR3263", "Dublin", "(DUB)", "Ryanair ", "B738", "(SP-RSW)", "1:38 PM", "0 hours, 8 minutes")
    ,(245, "Friday, Jan 13", "1:50 PM", "W61783", "Oslo", "(OSL)", "Wizz Air ", "A320", "(HA-LYS)", "2:02 PM", "0 hours, 12 minutes")
    ,(246, "Friday, Jan 13", "2:00 PM", "FR7271", "Vaxjo", "(VXO)", "Ryanair ", "B738", "(SP-RKQ)", "2:08 PM", "0 hours, 8 minutes")
    ,(247, "Friday, Jan 13", "2:20 PM", "LH1377", "Frankfurt", "(FRA)", "Lufthansa ", "CRJ9", "(D-ACNR)", "2:21 PM", "0 hours, 1 minutes")
)

$startRow = 216
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 12).Value = $row[10]
}